$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "57.559.67"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.56%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.451.41"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.24%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "507.87"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.86%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "133.14"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.91%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.995"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.48%  "
$ws.Range("E8").Value = "  -0.66%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.448.44"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.09%  "
$ws.Range("E10").Value = "  +0.34%  "
$ws.Range("E12").Value = "  -0.31%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.64"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.94%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.881.69"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.10%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "57.523.87"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.43%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "22.00"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.52%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000133"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.79%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.465.82"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.08%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.34"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.90%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.12"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.28%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "314.74"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.12%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.40"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.69%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.10%  "
$ws.Range("E24").Value = "  -2.38%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "65.47"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.05%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.995"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.41%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.542.71"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.54%  "
$ws.Range("E28").Value = "  -4.43%  "
$ws.Range("E29").Value = "  -1.81%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.60"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.65%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "174.06"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.24%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0₃0735"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.05%  "
$ws.Range("E33").Value = "  +0.01%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.19"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.24%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.13"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.05%  "
$ws.Range("E36").Value = "  +0.03%  "
$ws.Range("E37").Value = "  -0.22%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "17.96"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.10%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.24"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.73%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.86"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.33%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.819"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.92%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "36.55"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.14%  "
$ws.Range("E43").Value = "  +1.53%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "135.11"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +10.77%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.39"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.54%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.94"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.96%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "257.25"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.75%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.571"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.17%  "
$ws.Range("E49").Value = "  -0.22%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0495"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.09%  "
$ws.Range("E51").Value = "  +1.46%  "
